$wb = $excel.ActiveWorkbook

# Sheet ALC row 32 (hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1868.1
$ws.Range("I32").Value = 946.5
$ws.Range("J32").Value = 3250.5
$ws.Range("K32").Value = 946.5
$ws.Range("L32").Value = 3250.5
$ws.Range("M32").Value = -620.5
$ws.Range("N32").Value = -3902.5

# Sheet ALC row 64 (hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4041.2827
$ws.Range("I64").Value = 3890.476
$ws.Range("K64").Value = 3890.476
$ws.Range("M64").Value = -3642.476

# Sheet ALC row 67 (hunk 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4041.2827
$ws.Range("I67").Value = 3890.476
$ws.Range("K67").Value = 3890.476
$ws.Range("M67").Value = -3032.476

# Sheet ALC row 74 (hunk 3)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3600
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3600
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 3600
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -5472

# Sheet ALC row 76 (hunk 4)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3999.9524
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3999.9524
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 3999.9524
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -4629.9524

# Sheet ALC row 77 (hunk 5)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3600
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3600
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 18000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -27360

# Sheet ALC row 79 (hunk 6)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3999.9524
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3999.9524
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 3999.9524
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -6183.9524

# Sheet ALC row 119 (hunk 7)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H119").Value = 912.5
$ws.Range("J119").Value = 912.5
$ws.Range("L119").Value = 2737.5
$ws.Range("N119").Value = -12413.5

# Sheet ALC row 132 (hunk 8)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1539.2812
$ws.Range("I132").Value = 1666.5818
$ws.Range("J132").Value = 761.3333
$ws.Range("K132").Value = 4999.7454
$ws.Range("L132").Value = 2283.9999
$ws.Range("M132").Value = -2469.7454
$ws.Range("N132").Value = -7343.9999

# Sheet ALC row 137 (hunk 9)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 9655842
$ws.Range("I137").Value = 20834998
$ws.Range("J137").Value = 73707.21000000001
$ws.Range("K137").Value = 62504994
$ws.Range("L137").Value = 221121.63
$ws.Range("M137").Value = -62502444
$ws.Range("N137").Value = -226221.63

# Sheet ARM row 32 (hunk 10)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15154985
$ws.Range("I32").Value = 17243092
$ws.Range("J32").Value = 16201.75
$ws.Range("K32").Value = 17243092
$ws.Range("L32").Value = 16201.75
$ws.Range("M32").Value = -17242805
$ws.Range("N32").Value = -16775.75

# Sheet ARM row 61 (hunk 11)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3321.8076
$ws.Range("I61").Value = 1689.2307
$ws.Range("J61").Value = 4954.385
$ws.Range("K61").Value = 1689.2307
$ws.Range("L61").Value = 4954.385
$ws.Range("M61").Value = -1477.2307
$ws.Range("N61").Value = -5378.385

# Sheet ARM row 63 (hunk 12)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3715
$ws.Range("I63").Value = 2678
$ws.Range("J63").Value = 4337.2
$ws.Range("K63").Value = 2678
$ws.Range("L63").Value = 4337.2
$ws.Range("M63").Value = -1992
$ws.Range("N63").Value = -5709.2

# Sheet ARM row 66 (hunk 13)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3715
$ws.Range("I66").Value = 2678
$ws.Range("J66").Value = 4337.2
$ws.Range("K66").Value = 13390
$ws.Range("L66").Value = 21686
$ws.Range("M66").Value = -9958
$ws.Range("N66").Value = -28550

# Sheet ARM row 136 (hunk 14)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3321.8076
$ws.Range("I136").Value = 1689.2307
$ws.Range("J136").Value = 4954.385
$ws.Range("K136").Value = 5067.6921
$ws.Range("L136").Value = 14863.155
$ws.Range("M136").Value = -2517.6921
$ws.Range("N136").Value = -19963.155

# Sheet BSM row 2 (hunk 15)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# Sheet BSM row 86 (hunk 16)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1953.4
$ws.Range("I86").Value = 1753.9231
$ws.Range("K86").Value = 1753.9231
$ws.Range("M86").Value = -630.9231

# Sheet BSM row 89 (hunk 17)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1953.4
$ws.Range("I89").Value = 1753.9231
$ws.Range("K89").Value = 8769.6155
$ws.Range("M89").Value = -3153.6155

# Sheet BSM row 107 (hunk 18)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1838.1578
$ws.Range("I107").Value = 1987.5
$ws.Range("J107").Value = 1420
$ws.Range("K107").Value = 1987.5
$ws.Range("L107").Value = 1420
$ws.Range("M107").Value = -67.5
$ws.Range("N107").Value = -5260

# Sheet BSM row 115 (hunk 19)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H115").Value = 20750
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 20750
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 20750
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -23884

# Sheet CRP row 31 (hunk 20)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3108.12
$ws.Range("I31").Value = 2051.4707
$ws.Range("K31").Value = 2051.4707
$ws.Range("M31").Value = -1756.4707

# Sheet CRP row 34 (hunk 21)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3108.12
$ws.Range("I34").Value = 2051.4707
$ws.Range("K34").Value = 2051.4707
$ws.Range("M34").Value = -1849.4707

# Sheet CRP row 62 (hunk 22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2711.2856
$ws.Range("J62").Value = 2832.25
$ws.Range("L62").Value = 2832.25
$ws.Range("N62").Value = -4080.25

# Sheet CRP row 65 (hunk 23)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2711.2856
$ws.Range("J65").Value = 2832.25
$ws.Range("L65").Value = 14161.25
$ws.Range("N65").Value = -20401.25

# Sheet CUL row 5 (hunk 24)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1076.2632
$ws.Range("I5").Value = 1182.5555
$ws.Range("J5").Value = 980.6
$ws.Range("K5").Value = 3547.6665
$ws.Range("L5").Value = 2941.8
$ws.Range("M5").Value = -3435.6665
$ws.Range("N5").Value = -3165.8

# Sheet CUL row 56 (hunk 25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6666.6665
$ws.Range("I56").Value = 6666.6665
$ws.Range("K56").Value = 6666.6665
$ws.Range("M56").Value = -6136.6665

# Sheet CUL row 70 (hunk 26)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2667.3333
$ws.Range("I70").Value = 1282
$ws.Range("K70").Value = 3846
$ws.Range("M70").Value = -3531

# Sheet CUL row 73 (hunk 27)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 2667.3333
$ws.Range("I73").Value = 1282
$ws.Range("K73").Value = 3846
$ws.Range("M73").Value = -2754

# Sheet CUL row 135 (hunk 28)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1076.2632
$ws.Range("I135").Value = 1182.5555
$ws.Range("J135").Value = 980.6
$ws.Range("K135").Value = 10642.9995
$ws.Range("L135").Value = 8825.4
$ws.Range("M135").Value = -8107.9995
$ws.Range("N135").Value = -13895.4

# Sheet GSM row 70 (hunk 29)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8461.538
$ws.Range("I70").Value = 8461.538
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 8461.538
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -8191.538
$ws.Range("N70").ClearContents()

# Sheet GSM row 73 (hunk 30)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8461.538
$ws.Range("I73").Value = 8461.538
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 8461.538
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -7525.538
$ws.Range("N73").ClearContents()

# Sheet GSM row 80 (hunk 31)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2843.8667
$ws.Range("I80").Value = 2740.16
$ws.Range("J80").Value = 3362.4
$ws.Range("K80").Value = 2740.16
$ws.Range("L80").Value = 3362.4
$ws.Range("M80").Value = -1742.16
$ws.Range("N80").Value = -5358.4

# Sheet GSM row 83 (hunk 32)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2843.8667
$ws.Range("I83").Value = 2740.16
$ws.Range("J83").Value = 3362.4
$ws.Range("K83").Value = 13700.8
$ws.Range("L83").Value = 16812
$ws.Range("M83").Value = -8708.799999999999
$ws.Range("N83").Value = -26796

# Sheet GSM row 102 (hunk 33)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2531.4546
$ws.Range("I102").Value = 2427.3333
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2427.3333
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -805.3332999999998
$ws.Range("N102").Value = -6244

# Sheet GSM row 122 (hunk 34)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2320
$ws.Range("I122").Value = 2300
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 6900
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -4450
$ws.Range("N122").Value = -12400

# Sheet LTW row 40 (hunk 35)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6642.0835
$ws.Range("I40").Value = 10075
$ws.Range("J40").Value = 4925.625
$ws.Range("K40").Value = 10075
$ws.Range("L40").Value = 4925.625
$ws.Range("M40").Value = -9939
$ws.Range("N40").Value = -5197.625

# Sheet LTW row 74 (hunk 36)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 55000
$ws.Range("J74").Value = 55000
$ws.Range("L74").Value = 55000
$ws.Range("N74").Value = -56996

# Sheet LTW row 77 (hunk 37)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 55000
$ws.Range("J77").Value = 55000
$ws.Range("L77").Value = 165000
$ws.Range("N77").Value = -174984

# Sheet WVR row 81 (hunk 38)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1265.1666
$ws.Range("I81").Value = 897.75
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 1795.5
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -734.5
$ws.Range("N81").Value = -6122

# Sheet WVR row 84 (hunk 39)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1265.1666
$ws.Range("I84").Value = 897.75
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 8977.5
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -3673.5
$ws.Range("N84").Value = -30608

# Sheet WVR row 126 (hunk 40)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2633.1667
$ws.Range("I126").Value = 2866.3333
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 8598.999899999999
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -6128.999899999999
$ws.Range("N126").Value = -12140
